$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill D2:D5 with the same formula pattern used in D1 (SUM(A1,B1)/7+A1),
# applied relatively down the column, as a shared formula (D2 is master).
$ws.Range("D2:D5").Formula = "=SUM(A2,B2)/7+A2"
